$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 2 (the "decreased protein phosphorylation" record),
# shifting all rows below it up by one.
$ws.Rows("2").Delete()

# The row delete leaves one stale extra "Assigned By" hyperlink (carried
# over from the old last data row) sitting on the now-blank row 9. Rebuild
# the hyperlink list for the 7 remaining data rows (J2:J8) so it lines up
# 1:1 with the data again.
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 8; $r++) {
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 10), "mailto:robert-dodson@northwestern.edu") | Out-Null
}

# Adding hyperlinks re-applies the built-in "Hyperlink" cell style, which
# would otherwise register as a spurious duplicate style vs. the original
# formatting those cells already had. Copy the (still-original) formatting
# from the untouched trailing blank cell back over J2:J8 to keep the
# original style in place.
$ws.Range("J9").Copy()
$ws.Range("J2:J8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Select row 2 (now the "aberrant pseudopodium formation / DBS0351316" record)
$ws.Range("A2:XFD2").Select()
